$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.812.82"
$ws.Range("E2").Value = "  +4.16%  "
$ws.Range("D3").Value = "2.269.68"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.67"
$ws.Range("E5").Value = "  +4.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.06"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +3.87%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.79"
$ws.Range("E10").Value = "  +6.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.87"
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "2.621.30"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "2.276.58"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.764"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("D19").Value = "41.750.67"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.20"
$ws.Range("E20").Value = "  +8.16%  "
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.92"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.95"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.89"
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.92"
$ws.Range("E27").Value = "  +5.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.20"
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("E29").Value = "  +11.52%  "
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.25"
$ws.Range("E31").Value = "  +7.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.30"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.17"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0747"
$ws.Range("E35").Value = "  +4.57%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.15"
$ws.Range("E37").Value = "  +9.54%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.117"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").Value = "2.067.64"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.31"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.29"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("E47").Value = "  +5.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +6.86%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.67"
$ws.Range("E51").Value = "  +7.05%  "
